# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect a
# completed handback: the "Ready for handoff" status becomes
# "Handed back: in sync with en-US" on the Overview sheet, and each
# language sheet (zh-cn, de-de) gets its "Latest Target File" /
# "Latest Handback File" / "Latest Handback DateTime" columns filled
# in with the generated handback report info.

$wb = $excel.ActiveWorkbook

$targetFileName = "919668ce-0deb-44e5-82e5-c704e015e9a7.md"
$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cb8f71f7e727a64e003c6754c7d22b192a7a360f/e2e/919668ce-0deb-44e5-82e5-c704e015e9a7.md"

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: status columns (zh-cn / de-de) move from "Ready for
#     handoff" to "Handed back: in sync with en-US" ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus

$zhcnI2 = $zhcn.Range("I2")
$zhcn.Hyperlinks.Add($zhcnI2, $targetUrl, "", "", $targetFileName) | Out-Null
$zhcnI2.Style = "Hyperlink"

$zhcn.Range("J2").Value = "919668ce-0deb-44e5-82e5-c704e015e9a7.3df374e71a09553796dbe5c4c0a9960c02d4e655.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-27 19:05:45"

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus

$dedeI2 = $dede.Range("I2")
$dede.Hyperlinks.Add($dedeI2, $targetUrl, "", "", $targetFileName) | Out-Null
$dedeI2.Style = "Hyperlink"

$dede.Range("J2").Value = "919668ce-0deb-44e5-82e5-c704e015e9a7.3df374e71a09553796dbe5c4c0a9960c02d4e655.de-de.xlf"
$dede.Range("K2").Value = "2016-08-27 19:05:52"

# --- Column widths: widen the status/date columns that now hold longer
#     handback text ---
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$zhcn.Columns.Item(9).ColumnWidth = 40
$zhcn.Columns.Item(10).ColumnWidth = 40

$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(9).ColumnWidth = 40
$dede.Columns.Item(10).ColumnWidth = 40
